# Automatic tracker update ("Actualización automática del tracker")
#
# 1) Fill in the outcome ("resultado") and profit ("profit") for the last
#    previously-unresolved match (row 143).
# 2) Append three brand-new matches (rows 144-146) whose outcome is not
#    known yet, so their "resultado"/"profit" columns are left blank -
#    matching every other still-pending row already in the sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- 1) Resolve row 143 -----------------------------------------------
$ws.Range("G143").Value = "Acierto"
$ws.Range("H143").Value = 0.73

# --- 2) Append the new pending matches ----------------------------------
$newRows = @(
    @(14869898, "2025-10-12", "Sergey Fomin",   "Wishaya Trongcharoenchaikul", "Gana Wishaya Trongcharoenchaikul", 2.75),
    @(14869907, "2025-10-12", "Egor Agafonov",  "Mitsuki Wei Kang Leong",      "Gana Mitsuki Wei Kang Leong",      2.38),
    @(14869899, "2025-10-12", "Luca Castelnuovo","Duckhee Lee",                "Gana Duckhee Lee",                 2.75)
)

$startRow = 144
$endRow = $startRow + $newRows.Count - 1

# Force the "fecha" column to stay plain text (it is a yyyy-mm-dd *string*
# in this tracker, not a real date) and revert to the default style
# afterwards so no stray formatting is left on the cells.
$ws.Range("B$startRow`:B$endRow").NumberFormat = "@"

for ($i = 0; $i -lt $newRows.Count; $i++) {
    $r = $startRow + $i
    $row = $newRows[$i]
    $ws.Cells.Item($r, 1).Value = $row[0]
    $ws.Cells.Item($r, 2).Value = $row[1]
    $ws.Cells.Item($r, 3).Value = $row[2]
    $ws.Cells.Item($r, 4).Value = $row[3]
    $ws.Cells.Item($r, 5).Value = $row[4]
    $ws.Cells.Item($r, 6).Value = $row[5]
    # G (resultado) / H (profit) intentionally left blank - match isn't
    # settled yet, same as the other still-pending rows in the sheet.
}

$ws.Range("B$startRow`:B$endRow").Style = "Normal"
